$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 64 and 79 had placeholder 0 values in columns C and E:Z.
# These are cleared to blank (empty) cells, matching the consolidated
# ("balanços concatenados") sheet where those rows should have no data
# in those columns instead of numeric zeros.
$ws.Range("C64").ClearContents()
$ws.Range("E64:Z64").ClearContents()

$ws.Range("C79").ClearContents()
$ws.Range("E79:Z79").ClearContents()
